# Updates loading_percent values for Case_3_32 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 11.02949592498183
$ws.Range("C2").Value = 6.534646397226478
$ws.Range("D2").Value = 5.952329178425861
$ws.Range("F2").Value = 29.60148885121266
$ws.Range("G2").Value = 3.664825440854381
$ws.Range("I2").Value = 23.63300401273158
$ws.Range("K2").Value = 11.41412099186364
$ws.Range("N2").Value = 19.39942542237698

# Row 3
$ws.Range("B3").Value = 10.76789784190014
$ws.Range("C3").Value = 6.334895046498746
$ws.Range("D3").Value = 5.93506451983622
$ws.Range("F3").Value = 29.61492478565791
$ws.Range("G3").Value = 3.66714188255632
$ws.Range("I3").Value = 23.71487049766979
$ws.Range("K3").Value = 11.24090438113097
$ws.Range("N3").Value = 19.46779707327449

# Row 4
$ws.Range("B4").Value = 10.6067064918257
$ws.Range("C4").Value = 6.210928206341651
$ws.Range("D4").Value = 5.924314909497383
$ws.Range("F4").Value = 29.63145243576174
$ws.Range("G4").Value = 3.668639123213434
$ws.Range("I4").Value = 23.77005964855058
$ws.Range("K4").Value = 11.1359417467182
$ws.Range("N4").Value = 19.51170662294467

# Row 5
$ws.Range("B5").Value = 10.54098092910881
$ws.Range("C5").Value = 6.160167599802673
$ws.Range("D5").Value = 5.919897850493808
$ws.Range("F5").Value = 29.64026528884904
$ws.Range("G5").Value = 3.669268165634274
$ws.Range("I5").Value = 23.79378454918449
$ws.Range("K5").Value = 11.09357480768847
$ws.Range("N5").Value = 19.53008664718057

# Row 6
$ws.Range("B6").Value = 10.5300680237584
$ws.Range("C6").Value = 6.151726847544135
$ws.Range("D6").Value = 5.919162224725103
$ws.Range("F6").Value = 29.64185400694687
$ws.Range("G6").Value = 3.66937376125172
$ws.Range("I6").Value = 23.79779854546547
$ws.Range("K6").Value = 11.08656596679512
$ws.Range("N6").Value = 19.53316806364162

# Row 7
$ws.Range("B7").Value = 10.60582011051748
$ws.Range("C7").Value = 6.210244494511008
$ws.Range("D7").Value = 5.924255486034845
$ws.Range("F7").Value = 29.63156288252911
$ws.Range("G7").Value = 3.668647530073657
$ws.Range("I7").Value = 23.77037461494124
$ws.Range("K7").Value = 11.13536865322479
$ws.Range("N7").Value = 19.51195253055613

# Row 8
$ws.Range("B8").Value = 10.93947646089402
$ws.Range("C8").Value = 6.466097205693705
$ws.Range("D8").Value = 5.946407838227123
$ws.Range("F8").Value = 29.60440153572828
$ws.Range("G8").Value = 3.665608633999325
$ws.Range("I8").Value = 23.66020805878046
$ws.Range("K8").Value = 11.35414040687754
$ws.Range("N8").Value = 19.42260046054172

# Row 9
$ws.Range("B9").Value = 11.58487833181121
$ws.Range("C9").Value = 6.953761977664416
$ws.Range("D9").Value = 5.988624366722171
$ws.Range("F9").Value = 29.6169373720515
$ws.Range("G9").Value = 3.660241148267918
$ws.Range("I9").Value = 23.4833647423626
$ws.Range("K9").Value = 11.79174085446894
$ws.Range("N9").Value = 19.26262053550528

# Row 10
$ws.Range("B10").Value = 12.04809858951048
$ws.Range("C10").Value = 7.299089151912932
$ws.Range("D10").Value = 6.018840376822291
$ws.Range("F10").Value = 29.66634026391367
$ws.Range("G10").Value = 3.65665449381408
$ws.Range("I10").Value = 23.37751030100131
$ws.Range("K10").Value = 12.1152063093897
$ws.Range("N10").Value = 19.15427765092197

# Row 11
$ws.Range("B11").Value = 12.25538827372165
$ws.Range("C11").Value = 7.45257157253632
$ws.Range("D11").Value = 6.032401129036549
$ws.Range("F11").Value = 29.69753072400184
$ws.Range("G11").Value = 3.655099479519545
$ws.Range("I11").Value = 23.3346188280399
$ws.Range("K11").Value = 12.26210797595586
$ws.Range("N11").Value = 19.10696565136973

# Row 12
$ws.Range("B12").Value = 12.33330956635568
$ws.Range("C12").Value = 7.51011411855426
$ws.Range("D12").Value = 6.037508667566159
$ws.Range("F12").Value = 29.71059174006479
$ws.Range("G12").Value = 3.654521583694201
$ws.Range("I12").Value = 23.31913665353585
$ws.Range("K12").Value = 12.31764616975879
$ws.Range("N12").Value = 19.08933213816894

# Row 13
$ws.Range("B13").Value = 12.31655461707183
$ws.Range("C13").Value = 7.497747880217549
$ws.Range("D13").Value = 6.036409913365187
$ws.Range("F13").Value = 29.70772330140881
$ws.Range("G13").Value = 3.654645557602363
$ws.Range("I13").Value = 23.32243717028382
$ws.Range("K13").Value = 12.30568989854397
$ws.Range("N13").Value = 19.09311728648818

# Row 14
$ws.Range("B14").Value = 12.26181086971427
$ws.Range("C14").Value = 7.457317531277857
$ws.Range("D14").Value = 6.032821889057077
$ws.Range("F14").Value = 29.69858023833508
$ws.Range("G14").Value = 3.655051716464138
$ws.Range("I14").Value = 23.33332985278919
$ws.Range("K14").Value = 12.26667925488503
$ws.Range("N14").Value = 19.10550927705206

# Row 15
$ws.Range("B15").Value = 12.22820157608072
$ws.Range("C15").Value = 7.432475865318679
$ws.Range("D15").Value = 6.030620491364952
$ws.Range("F15").Value = 29.69314248564233
$ws.Range("G15").Value = 3.655301925396363
$ws.Range("I15").Value = 23.34010098954071
$ws.Range("K15").Value = 12.24277074852497
$ws.Range("N15").Value = 19.11313648225522

# Row 16
$ws.Range("B16").Value = 12.03447557524302
$ws.Range("C16").Value = 7.288980918146986
$ws.Range("D16").Value = 6.017950342427902
$ws.Range("F16").Value = 29.66447702432418
$ws.Range("G16").Value = 3.656757653636213
$ws.Range("I16").Value = 23.38041952397106
$ws.Range("K16").Value = 12.1055962662065
$ws.Range("N16").Value = 19.1574092059009

# Row 17
$ws.Range("B17").Value = 11.91469222190779
$ws.Range("C17").Value = 7.199983248496333
$ws.Range("D17").Value = 6.010129558649938
$ws.Range("F17").Value = 29.64912218348846
$ws.Range("G17").Value = 3.657670267392246
$ws.Range("I17").Value = 23.40650372053008
$ws.Range("K17").Value = 12.02134175515319
$ws.Range("N17").Value = 19.18507367345979

# Row 18
$ws.Range("B18").Value = 11.84547749187413
$ws.Range("C18").Value = 7.148457613953925
$ws.Range("D18").Value = 6.005613933310613
$ws.Range("F18").Value = 29.64111126368639
$ws.Range("G18").Value = 3.658202389445024
$ws.Range("I18").Value = 23.42200178145536
$ws.Range("K18").Value = 11.97286097453916
$ws.Range("N18").Value = 19.20117139673198

# Row 19
$ws.Range("B19").Value = 11.82199044785472
$ws.Range("C19").Value = 7.130956017366365
$ws.Range("D19").Value = 6.004082078504182
$ws.Range("F19").Value = 29.63853995929317
$ws.Range("G19").Value = 3.658383797020858
$ws.Range("I19").Value = 23.42733409764954
$ws.Range("K19").Value = 11.95644452258325
$ws.Range("N19").Value = 19.20665376967779

# Row 20
$ws.Range("B20").Value = 11.927476938069
$ws.Range("C20").Value = 7.209492478380159
$ws.Range("D20").Value = 6.010963894457856
$ws.Range("F20").Value = 29.65067180643274
$ws.Range("G20").Value = 3.657572372225121
$ws.Range("I20").Value = 23.4036757467198
$ws.Range("K20").Value = 12.03031324777074
$ws.Range("N20").Value = 19.18210951801472

# Row 21
$ws.Range("B21").Value = 12.27790665360286
$ws.Range("C21").Value = 7.469209029644295
$ws.Range("D21").Value = 6.03387653700046
$ws.Range("F21").Value = 29.70123189257642
$ws.Range("G21").Value = 3.654932120913046
$ws.Range("I21").Value = 23.33010975950367
$ws.Range("K21").Value = 12.27814050657259
$ws.Range("N21").Value = 19.1018617920199

# Row 22
$ws.Range("B22").Value = 12.50354886564696
$ws.Range("C22").Value = 7.635557091368154
$ws.Range("D22").Value = 6.048689763994513
$ws.Range("F22").Value = 29.74155830117429
$ws.Range("G22").Value = 3.653270388300629
$ws.Range("I22").Value = 23.28646081056665
$ws.Range("K22").Value = 12.43955956565449
$ws.Range("N22").Value = 19.05106137622593

# Row 23
$ws.Range("B23").Value = 12.38345429177897
$ws.Range("C23").Value = 7.547102293047216
$ws.Range("D23").Value = 6.040798797265746
$ws.Range("F23").Value = 29.71937056046982
$ws.Range("G23").Value = 3.654151464984674
$ws.Range("I23").Value = 23.30935063506469
$ws.Range("K23").Value = 12.35347524219442
$ws.Range("N23").Value = 19.07802432284502

# Row 24
$ws.Range("B24").Value = 11.92169805081342
$ws.Range("C24").Value = 7.205194471357035
$ws.Range("D24").Value = 6.010586751419498
$ws.Range("F24").Value = 29.64996867736639
$ws.Range("G24").Value = 3.65761660742468
$ws.Range("I24").Value = 23.40495271045352
$ws.Range("K24").Value = 12.02625736172737
$ws.Range("N24").Value = 19.18344901120749

# Row 25
$ws.Range("B25").Value = 11.4118242719515
$ws.Range("C25").Value = 6.823818575234708
$ws.Range("D25").Value = 5.977340318204705
$ws.Range("F25").Value = 29.60648808966737
$ws.Range("G25").Value = 3.661630248792819
$ws.Range("I25").Value = 23.52699093163784
$ws.Range("K25").Value = 11.67277991773348
$ws.Range("N25").Value = 19.30427783391555
